$d = $word.ActiveDocument

# The target content lives inside a legacy VML <w:pict> drawing (a v:group
# of v:rect/v:shape elements), which isn't reachable through the normal
# Word text object model (Find, Paragraphs, Shapes, ...). Work on the
# document's raw OOXML instead.
$xml = $d.Content.XML()

# Anchor: right after the "Hardware" textbox rect closes, and right before
# the existing <v:shapetype id="_x0000_t34"> definition used later by shape
# _x0000_s1047. We insert a (duplicate) shapetype definition together with
# a brand new connector shape (_x0000_s1035, rotation:180, no flip) in
# front of it, matching the upstream fix for 180 degree rotated shapes.
$closing = '</v:textbox>' + "`n" + '            </v:rect>' + "`n"
$anchor = $closing + '            <v:shapetype id="_x0000_t34"'

$idx = $xml.IndexOf($anchor)
if ($idx -lt 0) {
    throw "anchor not found in document XML"
}

$insertPos = $idx + $closing.Length

$insertion = '            <v:shapetype id="_x0000_t34" coordsize="21600,21600" o:spt="34" o:oned="t" adj="10800" path="m,l@0,0@0,21600,21600,21600e" filled="f">' + "`n" +
             '              <v:stroke joinstyle="miter"/>' + "`n" +
             '              <v:formulas>' + "`n" +
             '                <v:f eqn="val #0"/>' + "`n" +
             '              </v:formulas>' + "`n" +
             '              <v:path arrowok="t" fillok="f" o:connecttype="none"/>' + "`n" +
             '              <v:handles>' + "`n" +
             '                <v:h position="#0,center"/>' + "`n" +
             '              </v:handles>' + "`n" +
             '              <o:lock v:ext="edit" shapetype="t"/>' + "`n" +
             '            </v:shapetype>' + "`n" +
             '            <v:shape id="_x0000_s1035" type="#_x0000_t34" style="position:absolute;left:2956;top:291;width:1;height:495;rotation:180" o:connectortype="elbow" adj="-7776000,-486628,77954400">' + "`n" +
             '              <v:stroke startarrow="block" endarrow="block"/>' + "`n" +
             '            </v:shape>' + "`n"

$newXml = $xml.Substring(0, $insertPos) + $insertion + $xml.Substring($insertPos)

$d.Content.InsertXML($newXml)
